$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.337.68"
$ws.Range("E2").Value = "  +4.01%  "

$ws.Range("D3").Value = "1.732.21"
$ws.Range("E3").Value = "  +2.42%  "

$ws.Range("E4").Value = "  -0.07%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.33"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("E6").Value = "  -0.07%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -0.09%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.01"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +3.92%  "

$ws.Range("E9").Value = "  +1.74%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0637"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.36%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "1.977.33"
$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").Value = "1.732.79"
$ws.Range("E13").Value = "  +2.75%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.25"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("E15").Value = "  +1.74%  "

$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").Value = "28.332.76"
$ws.Range("E17").Value = "  +4.06%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "247.29"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("D19").Value = "0.0₃0754"
$ws.Range("E19").Value = "  +0.93%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.93"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -2.85%  "

$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").Value = "  -1.06%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.59"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("E26").Value = "  +2.35%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.73"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("E30").Value = "  +2.56%  "

$ws.Range("E31").Value = "  +2.27%  "

$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("D34").Value = "1.491.42"
$ws.Range("E34").Value = "  -5.77%  "

$ws.Range("E35").Value = "  -2.07%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.983"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +2.53%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  +0.59%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0177"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +1.38%  "

$ws.Range("E40").Value = "  +0.16%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.13"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("E42").Value = "  -0.14%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.66"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("D45").Value = "1.880.74"
$ws.Range("E45").Value = "  +2.18%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.800"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +1.69%  "

$ws.Range("E47").Value = "  +7.05%  "

$ws.Range("E48").Value = "  +4.36%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "90.51"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -0.97%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.19"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("E51").Value = "  -1.02%  "

